$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.371985244661392
$ws.Range("C2").Value = 1.100385049623747
$ws.Range("D2").Value = 3.542583507443938
$ws.Range("E2").Value = 1.034636339911349
$ws.Range("F2").Value = 3.524883908303043
$ws.Range("G2").Value = 0.628470476361786
$ws.Range("H2").Value = 3.347340434915449
$ws.Range("I2").Value = 0.9062051373079006
$ws.Range("J2").Value = 3.323262318674384
$ws.Range("K2").Value = 0.6290185522913129
